$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2055.45
$ws.Range("I15").Value = 2055.45
$ws.Range("K15").Value = 6166.349999999999
$ws.Range("M15").Value = -5997.349999999999

$ws.Range("H43").Value = 1377.6666
$ws.Range("I43").Value = 1498.6154
$ws.Range("J43").Value = 1234.7273
$ws.Range("K43").Value = 1498.6154
$ws.Range("L43").Value = 1234.7273
$ws.Range("M43").Value = -1429.6154
$ws.Range("N43").Value = -1372.7273

$ws.Range("H76").Value = 5069.4165
$ws.Range("I76").Value = 5416.4443
$ws.Range("J76").Value = 4028.3333
$ws.Range("K76").Value = 5416.4443
$ws.Range("L76").Value = 4028.3333
$ws.Range("M76").Value = -5101.4443
$ws.Range("N76").Value = -4658.3333

$ws.Range("H79").Value = 5069.4165
$ws.Range("I79").Value = 5416.4443
$ws.Range("J79").Value = 4028.3333
$ws.Range("K79").Value = 5416.4443
$ws.Range("L79").Value = 4028.3333
$ws.Range("M79").Value = -4324.4443
$ws.Range("N79").Value = -6212.3333

$ws.Range("H86").Value = 2127.4546
$ws.Range("I86").Value = 2062.875
$ws.Range("J86").Value = 2299.6667
$ws.Range("K86").Value = 2062.875
$ws.Range("L86").Value = 2299.6667
$ws.Range("M86").Value = -939.875
$ws.Range("N86").Value = -4545.6667

$ws.Range("H89").Value = 2127.4546
$ws.Range("I89").Value = 2062.875
$ws.Range("J89").Value = 2299.6667
$ws.Range("K89").Value = 10314.375
$ws.Range("L89").Value = 11498.3335
$ws.Range("M89").Value = -4698.375
$ws.Range("N89").Value = -22730.3335

$ws.Range("H92").Value = 1051.6471
$ws.Range("I92").Value = 1171.2858
$ws.Range("K92").Value = 1171.2858
$ws.Range("M92").Value = 76.71419999999989

$ws.Range("H113").Value = 3062.6924
$ws.Range("I113").Value = 2382.5
$ws.Range("K113").Value = 2382.5
$ws.Range("M113").Value = 871.5

$ws.Range("H132").Value = 2058.4082
$ws.Range("I132").Value = 1820
$ws.Range("J132").Value = 2988.2
$ws.Range("K132").Value = 5460
$ws.Range("L132").Value = 8964.599999999999
$ws.Range("M132").Value = -2930
$ws.Range("N132").Value = -14024.6

$ws.Range("H137").Value = 1143.0333
$ws.Range("I137").Value = 1122.2
$ws.Range("J137").Value = 1247.2
$ws.Range("K137").Value = 3366.6
$ws.Range("L137").Value = 3741.6
$ws.Range("M137").Value = -816.6000000000004
$ws.Range("N137").Value = -8841.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1668.4546
$ws.Range("I45").Value = 1478.625
$ws.Range("J45").Value = 2174.6667
$ws.Range("K45").Value = 1478.625
$ws.Range("L45").Value = 2174.6667
$ws.Range("M45").Value = -1101.625
$ws.Range("N45").Value = -2928.6667

$ws.Range("H74").Value = 903.59576
$ws.Range("I74").Value = 818.5135
$ws.Range("J74").Value = 1218.4
$ws.Range("K74").Value = 818.5135
$ws.Range("L74").Value = 1218.4
$ws.Range("M74").Value = 55.48649999999998
$ws.Range("N74").Value = -2966.4

$ws.Range("H77").Value = 903.59576
$ws.Range("I77").Value = 818.5135
$ws.Range("J77").Value = 1218.4
$ws.Range("K77").Value = 4092.5675
$ws.Range("L77").Value = 6092
$ws.Range("M77").Value = 275.4324999999999
$ws.Range("N77").Value = -14828

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 144543.58
$ws.Range("I86").Value = 1999.6
$ws.Range("J86").Value = 500903.5
$ws.Range("K86").Value = 1999.6
$ws.Range("L86").Value = 500903.5
$ws.Range("M86").Value = -876.5999999999999
$ws.Range("N86").Value = -503149.5

$ws.Range("H89").Value = 144543.58
$ws.Range("I89").Value = 1999.6
$ws.Range("J89").Value = 500903.5
$ws.Range("K89").Value = 9998
$ws.Range("L89").Value = 2504517.5
$ws.Range("M89").Value = -4382
$ws.Range("N89").Value = -2515749.5

$ws.Range("H134").Value = 2418.48
$ws.Range("I134").Value = 2233.9062
$ws.Range("K134").Value = 6701.7186
$ws.Range("M134").Value = -4166.7186

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2103
$ws.Range("I31").Value = 1146.6923
$ws.Range("J31").Value = 4175
$ws.Range("K31").Value = 1146.6923
$ws.Range("L31").Value = 4175
$ws.Range("M31").Value = -851.6922999999999
$ws.Range("N31").Value = -4765

$ws.Range("H34").Value = 2103
$ws.Range("I34").Value = 1146.6923
$ws.Range("J34").Value = 4175
$ws.Range("K34").Value = 1146.6923
$ws.Range("L34").Value = 4175
$ws.Range("M34").Value = -944.6922999999999
$ws.Range("N34").Value = -4579

$ws.Range("H105").Value = 1504
$ws.Range("I105").Value = 1504.4445
$ws.Range("J105").Value = 1500
$ws.Range("K105").Value = 1504.4445
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = 242.5554999999999
$ws.Range("N105").Value = -4994

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 669.6111
$ws.Range("I2").Value = 1195.8889
$ws.Range("J2").Value = 143.33333
$ws.Range("K2").Value = 7175.3334
$ws.Range("L2").Value = 859.9999799999999
$ws.Range("M2").Value = -7062.3334
$ws.Range("N2").Value = -1085.99998

$ws.Range("H131").Value = 1627.3287
$ws.Range("I131").Value = 2038.3334
$ws.Range("J131").Value = 1590.5223
$ws.Range("K131").Value = 6115.0002
$ws.Range("L131").Value = 4771.5669
$ws.Range("M131").Value = -1075.0002
$ws.Range("N131").Value = -14851.5669

$ws.Range("H132").Value = 1912.72
$ws.Range("J132").Value = 2168.2
$ws.Range("L132").Value = 19513.8
$ws.Range("N132").Value = -24573.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 12503473
$ws.Range("I122").Value = 11366315
$ws.Range("J122").Value = 14290436
$ws.Range("K122").Value = 34098945
$ws.Range("L122").Value = 42871308
$ws.Range("M122").Value = -34096495
$ws.Range("N122").Value = -42876208

$ws.Range("H132").Value = 5257.1333
$ws.Range("I132").Value = 5262.375
$ws.Range("K132").Value = 15787.125
$ws.Range("M132").Value = -13257.125

$ws.Range("H133").Value = 84660
$ws.Range("J133").Value = 84660
$ws.Range("L133").Value = 84660
$ws.Range("N133").Value = -89720

$ws.Range("H136").Value = 22956700
$ws.Range("I136").Value = 34484140
$ws.Range("K136").Value = 103452420
$ws.Range("M136").Value = -103449870

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1932.6666
$ws.Range("I96").Value = 1932.6666
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1932.6666
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -559.6666
$ws.Range("N96").ClearContents()

$ws.Range("H104").Value = 122597.8
$ws.Range("J104").Value = 122597.8
$ws.Range("L104").Value = 122597.8
$ws.Range("N104").Value = -129585.8

$ws.Range("H113").Value = 433.46155
$ws.Range("I113").Value = 496.1111
$ws.Range("J113").Value = 292.5
$ws.Range("K113").Value = 1488.3333
$ws.Range("L113").Value = 877.5
$ws.Range("M113").Value = 681.6667
$ws.Range("N113").Value = -5217.5

$ws.Range("H132").Value = 1601.0492
$ws.Range("I132").Value = 957.3488
$ws.Range("K132").Value = 2872.0464
$ws.Range("M132").Value = -342.0464000000002

Write-Output "Applied all market-data updates"
